# issue #5: add legislator_id, name, date into dataframe
#
# The original export only wrote each "表" (table) without the
# legislator's identifying columns. This adds three trailing columns -
# date / legislator_name / legislator_id - to the 股票 (stock) sheet,
# mirroring what the scraper now embeds in every table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "羅淑蕾"
$legislatorId   = 1638
$reportDate     = "2012-05-01"

$lastRow = 11  # rows 2..11 hold the stock entries (A2:G11 in the original sheet)

# --- header row -------------------------------------------------------
# Clone the formatting of the existing header cell (bold font, border,
# centered alignment) onto the three new header cells so they render
# identically to B1:G1, then fill in the header text.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- data rows ----------------------------------------------------------
# Clone the plain data-cell formatting (no border, normal font) from an
# existing data cell onto the new columns for every data row.
$ws.Range("C2").Copy() | Out-Null
$ws.Range(("H2:J" + $lastRow)).PasteSpecial(-4122) | Out-Null

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value  = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}

# Force column H to text so the yyyy-mm-dd string is stored as a literal
# shared string rather than being auto-converted to a date serial number.
$ws.Range(("H2:H" + $lastRow)).NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $reportDate
}
